$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.990.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.906.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.11"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4634"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4064"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.82"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08002"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.67"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.909.07"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.939"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.098"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.97"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001033"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.47"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.982.19"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.462"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.247"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.133.52"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.15"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.71"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.103"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.396"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.05"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9800"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09379"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.418"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.297"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06078"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02228"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.399"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.165"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5813"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9994"
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1825"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.14"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.259"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.348"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +14.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.15"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5498"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.902"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07026"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.51"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +23.10%  "
